$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 3.06333513865784
$ws.Range("C1").Value = 4.969480335136774
$ws.Range("D1").Value = 6.903183631835497
$ws.Range("E1").Value = 8.921377202387516
$ws.Range("F1").Value = 10.87168641311016
$ws.Range("G1").Value = 12.79104239245185
$ws.Range("H1").Value = 14.80832674159909
$ws.Range("I1").Value = 16.62613124335098
$ws.Range("J1").Value = 18.73540606967725
$ws.Range("K1").Value = 20.8086201974007
